$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ayat = 'Surah Al Nisa, 14 - 17'
$tags = 'Colors of Quran, Punishment in Islam, Amputation of hand, Hanging the killer, Life in prison, Councelling the disbeliever'
$content = 'h1: Shocking misconception of punishments in Islam
quote: Those who commit<span class=lavendar> unlawful sexual intercourse </span>of your women - bring against them <span class=lavendar>four [witnesses] </span>from among you. And if they testify, <span class=lavendar>confine the guilty women to houses</span> until death takes them<span class=lavendar>or Allah ordains for them [another] way</span>. And the two who commit it among you, dishonor them both. But <span class=lavendar>if they repent and correct themselves, leave them alone</span>. Indeed, Allah is ever Accepting of repentance and Merciful. <br> - Surah Al Nisa verse 15 - 16
p: Quran tells of punishing the people committing unlawful sexual intercourse in a very precise manner. Each punishment is thoroughly elaborated and self sufficient. It also focuses on the leniency it offers to its offenders. Quran maintains the sanity of this world and does not enforce any <B>‘difficult to implement’</b> order. 
p: Above verse has following parts for punishing a fornicator:-
p.b-left: <b>1. Bring upon 4 witnesses. </b> If accuser is unable to bring upon these 4 witnesses, the punishment mentioned next, is inapplicable. Accused is free.
p.b-left: <b>2. Confine the guilty women. </b>It states of confinement to a closed house. She still can live inside this house, teach our kids, eat, live and does not undergo social victimisation.
p.b-left: <b>3. Allah ordains for them another way. </b>After the house arrest, it is also possible they do not die in there. They might are freed from the prison and they go out, if Allah wills. 
p: I believe, it is not possible to keep someone under house arrest for the entire length of his / her life. Circumstances and scenarios that bound us, usually always create a soft corner for the victim. Prolonged confinements to houses is only possible after it gets difficult for the Judiciary to work independently. As there gets more people involved, like for example in democracy, the law becomes a complicated mesh of papers and rules. As soon as the judge pours in ‘human’ part of himself, the offenders start becoming productive to society. They get leniency to learn skills in their confined places which leads them to becoming a contributing factor towards society.
p.note: These practices are happening in <b>Norway at Helden prison</b>. Offenders are kept in a confined school where they are taught useful skills of wiring, mechanics, programming and if they can’t read or write, they are taught to read and write. <a href=https://showsnob.com/2019/04/13/inside-the-worlds-toughest-prisons-season-3-finale/ target=_blank>Link to post is here.</a>
p.b-left: <b>4. If they repent and correct themselves, leave them alone. </b>This part is mostly missed out in all punishments in Fundamentalist Islam. If the offender in prison comes to repentance, free him.
h3: What confuses in Quran
quote: And yet the more you add, the more seems to go missing. <br> - Jew Writer Lesley Hazlton <a href=https://www.ted.com/talks/lesley_hazleton_the_doubt_essential_to_faith target=_blank>talking about Quran.</a>
p: As we starting adding meanings and explanations, we start getting closer to becoming a <b>fundamentalist</b>. The punishments’ color fades out. Kind and soft meanings laid out in Quran starts diminishing. As these are executed out in purity, the collapsing society starts getting back on its own feet.
p: Quran does not talk of murdering the infidel, it talks of counselling, embracing and inspiring the infidel. It talks of healing the broken fabrics around us.
h3: Conclusion
p: The punishments of ‘amputating’ the hand of a thief looks only severe because we do not know what it feels to sit in prison for an entire year in prison. 
p.b-left: 1 year in prison is greater than 1 less hand. 
p.b-left: 25 years in solitary confinements is greater than 1 x death. <br>The death does not impact any body for 25 long years. <br>The pain of losing someone ultimately is healed. <br>However, the pain of feeding your child in prison for 25 years is heavier.
p: Above all, if there is a way to avoid death, amputation or confinement of criminal, do it. Do not cut a thief’s hand, if he repents. Do not confine a fornicator, if he repents. Do not enforce the highlightened version of Quran. <b>Enforce the Quran in totality, that it does not loses its colors and dulls down our religion.</b>
p.note: Guys, if you like this project. Please follow this project''s page on twitter. <a href="https://twitter.com/zakatlists">Click here to go to the twitter page</a>. Subscribe below to receive updates.'

$ws.Range("A44").Value = 43
$ws.Range("B44").Value = 43873
$ws.Range("B44").NumberFormat = $ws.Range("B43").NumberFormat
$ws.Range("C44").Value = $ayat
$ws.Range("F44").Value = $tags
$ws.Range("D44").Value = $content
$ws.Range("E44").Value = "Qasim Ali"

$ws.Rows.Item(44).RowHeight = 409.6

$excel.ActiveWindow.ScrollRow = 44
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D44").Select()
